# Added alternate jobs plot.
#
# Starting point: workbook has tabs
#   Avg Demand Scenario, Blade, Nacelle, Tower, Array cable, Export cable,
#   Monopile, WTIV
#
# Target: three new "factory" tabs (Jacket, Semisubmersible,
# Offshore substation) are inserted between Monopile and WTIV, built by
# duplicating the existing WTIV factory-spec template, and the WTIV tab
# itself is populated with new figures. The "Avg Demand Scenario" roadmap
# sheet gets four new rows for the new facilities, and the Monopile sheet's
# selection/view is tidied up.

$wb = $excel.ActiveWorkbook

$missing = [System.Reflection.Missing]::Value

# ---------------------------------------------------------------------
# 1. Build the three new factory-spec sheets by copying the WTIV template
#    (it already has the right layout/labels/number formats) and inserting
#    each copy immediately before the WTIV tab so the final tab order is
#    ... Monopile, Jacket, Semisubmersible, Offshore substation, WTIV.
# ---------------------------------------------------------------------
function New-SheetBeforeWtiv([string]$newName) {
    $wtiv = $wb.Worksheets.Item("WTIV")
    $idx = $wtiv.Index
    $wtiv.Copy($wtiv, $missing)
    $newSheet = $wb.Worksheets.Item($idx)
    $newSheet.Name = $newName
    return $newSheet
}

$jacket = New-SheetBeforeWtiv("Jacket")
$semisub = New-SheetBeforeWtiv("Semisubmersible")
$offshore = New-SheetBeforeWtiv("Offshore substation")

# ---------------------------------------------------------------------
# 2. "Jacket" factory specifications.
# ---------------------------------------------------------------------
$jacket.Range("B2").Value = 50
$jacket.Range("C2").Value = "jackets/year"
$jacket.Range("B3").Value = 50
$jacket.Range("B12").Value = 425
$jacket.Range("B13").Value = 0
$jacket.Range("B14").Value = 7.1
$jacket.Range("B15").Value = 2.4
$jacket.Range("B16").Value = 4.7
$jacket.Range("B17").Formula = "=100-SUM(B14:B16)-B18"
$jacket.Range("B18").Value = 11.8
$jacket.Range("A1:D18").Select()

# ---------------------------------------------------------------------
# 3. "Semisubmersible" is left as a straight duplicate of the original
#    WTIV factory spec (values unchanged).
# ---------------------------------------------------------------------
$semisub.Range("B4").Select()

# ---------------------------------------------------------------------
# 4. "Offshore substation" factory specifications.
# ---------------------------------------------------------------------
$offshore.Range("B2").Value = 1
$offshore.Range("C2").Value = "substation/year"
$offshore.Range("B3").Value = 10
$offshore.Range("B4").Value = 0
$offshore.Range("B12").Value = 100
$offshore.Range("B13").Value = 0
$offshore.Range("B14").Value = 7.1
$offshore.Range("B15").Value = 2.4
$offshore.Range("B16").Value = 4.7
$offshore.Range("B17").Formula = "=100-SUM(B14:B16)-B18"
$offshore.Range("B18").Value = 11.8
$offshore.Range("E13").Select()

# ---------------------------------------------------------------------
# 5. The "WTIV" tab itself gets replaced with new figures; the lower
#    workforce breakdown rows are left blank on this tab.
# ---------------------------------------------------------------------
$wtivFinal = $wb.Worksheets.Item("WTIV")
$wtivFinal.Range("B2").Value = 50
$wtivFinal.Range("C2").Value = "semisubs/year"
$wtivFinal.Range("B3").Value = 100
$wtivFinal.Range("B4").Value = 3
$wtivFinal.Range("B12").Value = 200
$wtivFinal.Range("B13").ClearContents()
$wtivFinal.Range("B14").ClearContents()
$wtivFinal.Range("B15").ClearContents()
$wtivFinal.Range("B16").ClearContents()
$wtivFinal.Range("B17").ClearContents()
$wtivFinal.Range("B18").ClearContents()
$wtivFinal.Range("N13").Select()

# ---------------------------------------------------------------------
# 6. "Avg Demand Scenario" roadmap sheet: four new facility rows.
# ---------------------------------------------------------------------
$roadmap = $wb.Worksheets.Item("Avg Demand Scenario")
$roadmap.Range("A14").Value = "Jacket 1"
$roadmap.Range("B14").Value = 2024
$roadmap.Range("A15").Value = "Semisubmersible 1"
$roadmap.Range("B15").Value = 2028
$roadmap.Range("A16").Value = "Semisubmersible 2"
$roadmap.Range("B16").Value = 2029
$roadmap.Range("A17").Value = "Semisubmersible 3"
$roadmap.Range("B17").Value = 2030
$roadmap.Columns.Item(1).ColumnWidth = 19.18
$roadmap.Range("B18").Select()

# ---------------------------------------------------------------------
# 7. "Monopile" sheet: tidy up the view/selection.
# ---------------------------------------------------------------------
$monopile = $wb.Worksheets.Item("Monopile")
$monopile.Range("A1:XFD1048576").Select()
$excel.ActiveWindow.ScrollRow = 1

# ---------------------------------------------------------------------
# 8. Leave the final active sheet as the originally-selected one.
# ---------------------------------------------------------------------
$roadmap.Activate()
